$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new "placements" account line ---

# Column A: label, formatted like the other account-name cells in column A
# (copy the visual format used by A4 before writing the new text)
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "placements"

# Column B: plain value, default formatting
$ws.Range("B7").Value = 0

# Columns C:D, formatted like the other early-year cells (same look as C4/B4)
$ws.Range("C4").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

# Columns E:S, formatted like the other populated-year numeric cells (M2 style)
$ws.Range("M2").Copy()
$ws.Range("E7:S7").PasteSpecial(-4122)

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 36
$ws.Range("G7").Value = -50
$ws.Range("H7").Value = 60
$ws.Range("I7").Value = 94
$ws.Range("J7").Value = 149
$ws.Range("K7").Value = 232
$ws.Range("L7").Value = 339
$ws.Range("M7").Value = 315
$ws.Range("N7").Value = 298
$ws.Range("O7").Value = 422
$ws.Range("P7").Value = 412
$ws.Range("Q7").Value = 1394
$ws.Range("R7").Value = 464
$ws.Range("S7").Value = 769

# --- Move the active selection, matching the saved cursor position ---
$ws.Range("E11").Select() | Out-Null
